$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the stock id on rows 29 and 34: numeric 8289 -> text "8289_" ---
$ws.Range("B29").Value = "8289_"
$ws.Range("B34").Value = "8289_"

# --- Append the 2021/11/26 trading rows (52-75) ---
$dateFmt = "m""月""d""日"""

$newRows = @(
    @(52, 44526, 3532, "long",  12,   229),
    @(53, 44526, 6167, "long",  200,  13.05),
    @(54, 44526, 3444, "long",  43,   68.5),
    @(55, 44526, 3046, "long",  75,   39.3),
    @(56, 44526, 3653, "long",  7,    405),
    @(57, 44526, 6182, "long",  37,   78.9),
    @(58, 44526, 6457, "long",  13,   226.5),
    @(59, 44526, 3037, "long",  14,   211.5),
    @(60, 44526, 3169, "long",  11,   264.5),
    @(61, 44526, 3672, "long",  77,   38.4),
    @(62, 44526, 2374, "long",  125,  23.5),
    @(63, 44526, 6265, "long",  170,  17.4),
    @(64, 44526, 2383, "long",  11,   265.5),
    @(65, 44526, 3551, "long",  78,   75),
    @(66, 44526, 2340, "short", -57,  56.7),
    @(67, 44526, 6104, "short", -18,  178),
    @(68, 44526, 6138, "short", -15,  254),
    @(69, 44526, 3016, "short", -21,  141.5),
    @(70, 44526, 8069, "short", -25,  119),
    @(71, 44526, 8289, "short", -75,  39.7),
    @(72, 44526, 1712, "short", -107, 28.5),
    @(73, 44526, 3588, "short", -35,  152.5),
    @(74, 44526, 5351, "short", -63,  87),
    @(75, 44526, 2484, "short", -137, 41.1)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("A$r").NumberFormat = $dateFmt
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# --- Move the selection to reflect where the user left off editing ---
$ws.Range("F29").Select()
